# Update the "Stretch User Deliverables" bullet list: the two bullets about
# the card-flip / card-matching interaction are replaced with a single
# bullet describing the new film-trailer deliverable.

$d = $word.ActiveDocument

# Locate the two paragraphs to merge/rewrite via their distinctive text.
$oldFirstText  = "As a user, when I click on a card, it should flip to reveal a movie poster, with detail information about the movie (title, year, synopsis, cast)"
$oldSecondText = "As a user, if I click a second card and it matches, the cards will stay revealed (flipped over), if they do not match, they will flip back over."

$firstPara = $null
$secondPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    # Paragraph.Range.Text includes the trailing paragraph-mark character(s)
    # (CR / cell-mark / etc.), so trim those off before comparing.
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7, [char]11)
    if ($ptext -eq $oldFirstText) {
        $firstPara = $i
    }
    elseif ($ptext -eq $oldSecondText) {
        $secondPara = $i
    }
}

# Rewrite the first paragraph's run text in place (still its own paragraph).
$p1 = $d.Paragraphs($firstPara)
$p1.Range.Text = "As a user, "

# Rewrite the second paragraph's run text in place.
$p2 = $d.Paragraphs($secondPara)
$p2.Range.Text = "I will be able to see a trailer for the movie I searched for. The trailer will show inside the “theater” which will load next to the movie information. "

# Merge the two paragraphs into one by deleting the paragraph mark that
# ends the (now-shortened) first paragraph.
$p1 = $d.Paragraphs($firstPara)
$markStart = $p1.Range.End - 1
$d.Range($markStart, $markStart + 1).Delete()

Write-Output "Result: $($d.Paragraphs($firstPara).Range.Text)"
